$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D38").Value = "电话"
$ws.Range("D39").Value = "公司地址"
$ws.Range("D40").Value = "隐私权政策"
$ws.Range("D41").Value = "使用条款"
$ws.Range("D42").Value = "帮助中心"
$ws.Range("D43").Value = "© VIXI 公司名称，2022。版权所有"

$ws.Range("D38:D43").Font.Underline = $true

$ws.Range("G41").Select()
